$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TM")

# Row 7
$ws.Range("F7").Value = 702.1
$ws.Range("G7").Value = 711.9
$ws.Range("H7").Value = 699.4
$ws.Range("I7").Value = 708.85
$ws.Range("J7").Value = 705.25

# Row 9
$ws.Range("G9").Value = 705.25
$ws.Range("H9").Value = 696.25
$ws.Range("I9").Value = 700.4

# Row 10
$ws.Range("G10").Value = 706.9
$ws.Range("H10").Value = 699.55
$ws.Range("I10").Value = 706.9

# Row 11
$ws.Range("G11").Value = 711.4
$ws.Range("H11").Value = 706
$ws.Range("I11").Value = 708.6

# Row 12
$ws.Range("G12").Value = 710.2
$ws.Range("H12").Value = 705.7
$ws.Range("I12").Value = 710.2

# Row 13
$ws.Range("G13").Value = 710.75
$ws.Range("H13").Value = 707.25
$ws.Range("I13").Value = 708.95

# Row 14
$ws.Range("G14").Value = 711.4
$ws.Range("H14").Value = 708.1
$ws.Range("I14").Value = 710.25

# Row 15
$ws.Range("G15").Value = 711.7
$ws.Range("H15").Value = 710
$ws.Range("I15").Value = 711.45

# Row 16
$ws.Range("G16").Value = 711.9
$ws.Range("H16").Value = 710.5
$ws.Range("I16").Value = 711.2

# Row 17
$ws.Range("G17").Value = 711.5
$ws.Range("H17").Value = 709.1
$ws.Range("I17").Value = 710.8

# Row 18
$ws.Range("G18").Value = 711.7
$ws.Range("H18").Value = 707.15
$ws.Range("I18").Value = 708.35

# Row 19
$ws.Range("G19").Value = 708.8
$ws.Range("H19").Value = 704
$ws.Range("I19").Value = 708.4

# Row 20
$ws.Range("G20").Value = 710
$ws.Range("H20").Value = 707.8
$ws.Range("I20").Value = 707.95

# Row 21
$ws.Range("G21").Value = 710.75
$ws.Range("H21").Value = 707.3
$ws.Range("I21").Value = 709.95
